$d = $word.ActiveDocument

# 1. Update "Date of Quote:" text from 14-Sept-17 to 17-Sept-17 (appears twice:
#    DrawingML choice + VML fallback duplicate content)
$d.Content.Find.Execute(" 14-Sept-17", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " 17-Sept-17", 2)

# 2. Clear "project project description" text (appears twice)
$d.Content.Find.Execute("project project description", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 3. Clear "scope of work alenka" text (appears twice)
$d.Content.Find.Execute("scope of work alenka", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 4. Update Grand Total With Tax amount
$d.Content.Find.Execute("Grand Total With Tax: `$2.50", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Grand Total With Tax: `$0.00", 2)

# 5. Update header date from 13-Sept-17 to 17-Sept-17
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("Date: 13-Sept-17", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "Date: 17-Sept-17", 2)
    }
}

# 6. Remove the two demolition rows from the relevant table, leaving it empty
foreach ($tbl in $d.Tables) {
    $hasDemo = $false
    foreach ($r in $tbl.Range.Rows) {
        if ($tbl.Range.Text -like "*Demolition*") {
            $hasDemo = $true
        }
    }
}

foreach ($tbl in $d.Tables) {
    if ($tbl.Range.Text -like "*Demolition*") {
        while ($tbl.Rows.Count -gt 0) {
            $tbl.Rows(1).Delete()
        }
    }
}
